$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers (English)
$ws.Range("L1").Value = "itemsDiscount"
$ws.Range("M1").Value = "netTotal"
$ws.Range("N1").Value = "total"
$ws.Range("O1").Value = "T1(V009)"
$ws.Range("P1").Value = "T1(V001)"
$ws.Range("Q1").Value = "T2(Tbl01)"

# Row 2 headers (Arabic)
$ws.Range("L2").Value = "خصم الأصناف"
$ws.Range("M2").Value = "الصافى (بعد الخصم قبل الضريبة)"
$ws.Range("N2").Value = "الأجمالى الكلى"
$ws.Range("O2").Value = "T1(V009)"
$ws.Range("P2").Value = "T1(V001)"
$ws.Range("Q2").Value = "T2(Tbl01)"

# Row 3 data
$ws.Range("A3").Value = 500
$ws.Range("C3").Value = 30
$ws.Range("E3").Value = 4690
$ws.Range("F3").Value = "EG-237791390-QTC1012"
$ws.Range("H3").Value = 11
$ws.Range("K3").Value = 1100
$ws.Range("L3").Value = 100
$ws.Range("M3").Value = 1000
$ws.Range("N3").Value = 1145
$ws.Range("O3").Value = 140
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 5
